# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the row for
# 8a29e634-f382-48fb-80ee-1050e09a6216.md now sorts first (still "In
# Translation"), and the row for 3a646ba2-d225-41c5-abdc-7841eb2c7d46.md
# moves to second place and is updated to "Ready for handoff" (new
# handoff timestamps / mt priority). This touches the Overview sheet and
# the two per-locale sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6603960e2e11ef9397ffc025834a6424d3d739ab/e2e/"
$fileA = "8a29e634-f382-48fb-80ee-1050e09a6216.md"
$fileB = "3a646ba2-d225-41c5-abdc-7841eb2c7d46.md"

# ---------------------------------------------------------------------
# Overview sheet: row 2 -> fileA (unchanged status), row 3 -> fileB (now
# "Ready for handoff" with a later generate date).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $fileA
$wsOverview.Range("A3").Value = $fileB

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 06:14:46"

# Hyperlinks on column B must follow the swapped file names. This shim's
# Range.Hyperlinks collection mutates the whole sheet, so drop every
# hyperlink on the sheet and re-add the two, in order, to land on the
# same rIds as before.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $baseUrl + $fileA, "", "", "e2e\" + $fileA)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl + $fileB, "", "", "e2e\" + $fileB)

# Status column widened to fit "Ready for handoff".
$wsOverview.Range("E1:F1").ColumnWidth = 16.382520825846367

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): row 2 -> fileA, row 3 -> fileB
# (fileB's row also gets the new status / priority / handoff datetime).
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; Ext = "zh-cn.xlf"; H2 = "2016-08-19 06:14:10"; H3 = "2016-08-19 06:14:42" },
    @{ Name = "de-de"; Ext = "de-de.xlf"; H2 = "2016-08-19 06:14:15"; H3 = "2016-08-19 06:14:46" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    $ws.Range("A2").Value = $fileA
    $ws.Range("A3").Value = $fileB

    $ws.Range("G2").Value = "8a29e634-f382-48fb-80ee-1050e09a6216.a7abf00bd63f761bc7f3afd6fbf84767f6f991f7." + $locale.Ext
    $ws.Range("G3").Value = "3a646ba2-d225-41c5-abdc-7841eb2c7d46.243964ea6e64e0531922670ed4c4c6d3d65deb77." + $locale.Ext

    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("E3").Value = "mt"
    $ws.Range("H3").Value = $locale.H3

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $fileA, "", "", $fileA)
    $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $fileB, "", "", $fileB)

    $ws.Range("C1").ColumnWidth = 16.382520825846367
}
